$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 99.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 10755.9375
$ws.Range("J64").Value = 11936.77
$ws.Range("L64").Value = 11936.77
$ws.Range("N64").Value = -12432.77
$ws.Range("H67").Value = 10755.9375
$ws.Range("J67").Value = 11936.77
$ws.Range("L67").Value = 11936.77
$ws.Range("N67").Value = -13652.77
$ws.Range("H70").Value = 3561.111
$ws.Range("I70").Value = 1200
$ws.Range("K70").Value = 3600
$ws.Range("M70").Value = -3330
$ws.Range("H73").Value = 3561.111
$ws.Range("I73").Value = 1200
$ws.Range("K73").Value = 3600
$ws.Range("M73").Value = -2664
$ws.Range("H137").Value = 4459.6
$ws.Range("I137").Value = 2749.25
$ws.Range("K137").Value = 8247.75
$ws.Range("M137").Value = -5697.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5030.115
$ws.Range("I32").Value = 5474.413
$ws.Range("K32").Value = 5474.413
$ws.Range("M32").Value = -5187.413
$ws.Range("H61").Value = 16186.111
$ws.Range("I61").Value = 4400
$ws.Range("J61").Value = 19553.572
$ws.Range("K61").Value = 4400
$ws.Range("L61").Value = 19553.572
$ws.Range("M61").Value = -4188
$ws.Range("N61").Value = -19977.572
$ws.Range("H74").Value = 3489.25
$ws.Range("I74").Value = 3596.7856
$ws.Range("J74").Value = 3238.3333
$ws.Range("K74").Value = 3596.7856
$ws.Range("L74").Value = 3238.3333
$ws.Range("M74").Value = -2722.7856
$ws.Range("N74").Value = -4986.3333
$ws.Range("H77").Value = 3489.25
$ws.Range("I77").Value = 3596.7856
$ws.Range("J77").Value = 3238.3333
$ws.Range("K77").Value = 17983.928
$ws.Range("L77").Value = 16191.6665
$ws.Range("M77").Value = -13615.928
$ws.Range("N77").Value = -24927.6665
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 1418.2
$ws.Range("I132").Value = 1376.4286
$ws.Range("K132").Value = 4129.2858
$ws.Range("M132").Value = -1599.2858
$ws.Range("H136").Value = 16186.111
$ws.Range("I136").Value = 4400
$ws.Range("J136").Value = 19553.572
$ws.Range("K136").Value = 13200
$ws.Range("L136").Value = 58660.716
$ws.Range("M136").Value = -10650
$ws.Range("N136").Value = -63760.716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29472442
$ws.Range("I86").Value = 100201580
$ws.Range("J86").Value = 1966.0834
$ws.Range("K86").Value = 100201580
$ws.Range("L86").Value = 1966.0834
$ws.Range("M86").Value = -100200457
$ws.Range("N86").Value = -4212.0834
$ws.Range("H89").Value = 29472442
$ws.Range("I89").Value = 100201580
$ws.Range("J89").Value = 1966.0834
$ws.Range("K89").Value = 501007900
$ws.Range("L89").Value = 9830.416999999999
$ws.Range("M89").Value = -501002284
$ws.Range("N89").Value = -21062.417
$ws.Range("H94").Value = 7442.857
$ws.Range("I94").Value = 1050
$ws.Range("K94").Value = 1050
$ws.Range("M94").Value = -599
$ws.Range("H105").Value = 2848.7856
$ws.Range("I105").Value = 1969.1428
$ws.Range("J105").Value = 3728.4285
$ws.Range("K105").Value = 1969.1428
$ws.Range("L105").Value = 3728.4285
$ws.Range("M105").Value = -222.1428000000001
$ws.Range("N105").Value = -7222.4285
$ws.Range("H134").Value = 1881.7727
$ws.Range("I134").Value = 1500.3334
$ws.Range("K134").Value = 4501.0002
$ws.Range("M134").Value = -1966.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 444.6
$ws.Range("I22").Value = 443
$ws.Range("J22").Value = 448.33334
$ws.Range("K22").Value = 443
$ws.Range("L22").Value = 448.33334
$ws.Range("M22").Value = -93
$ws.Range("N22").Value = -1148.33334
$ws.Range("H31").Value = 3461.182
$ws.Range("I31").Value = 2562.7273
$ws.Range("J31").Value = 4359.636
$ws.Range("K31").Value = 2562.7273
$ws.Range("L31").Value = 4359.636
$ws.Range("M31").Value = -2267.7273
$ws.Range("N31").Value = -4949.636
$ws.Range("H34").Value = 3461.182
$ws.Range("I34").Value = 2562.7273
$ws.Range("J34").Value = 4359.636
$ws.Range("K34").Value = 2562.7273
$ws.Range("L34").Value = 4359.636
$ws.Range("M34").Value = -2360.7273
$ws.Range("N34").Value = -4763.636
$ws.Range("H58").Value = 2682.6667
$ws.Range("I58").Value = 1780.5
$ws.Range("K58").Value = 1780.5
$ws.Range("M58").Value = -1577.5
$ws.Range("H62").Value = 71434560
$ws.Range("I62").Value = 7579.8
$ws.Range("K62").Value = 7579.8
$ws.Range("M62").Value = -6955.8
$ws.Range("H65").Value = 71434560
$ws.Range("I65").Value = 7579.8
$ws.Range("K65").Value = 37899
$ws.Range("M65").Value = -34779
$ws.Range("H74").Value = 39710
$ws.Range("H77").Value = 39710
$ws.Range("H122").Value = 3776.353
$ws.Range("I122").Value = 2908.9092
$ws.Range("K122").Value = 8726.7276
$ws.Range("M122").Value = -6276.7276
$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -19910
$ws.Range("H134").Value = 5327.2666
$ws.Range("I134").Value = 4732.8
$ws.Range("K134").Value = 14198.4
$ws.Range("M134").Value = -11663.4
$ws.Range("H136").Value = 2682.6667
$ws.Range("I136").Value = 1780.5
$ws.Range("K136").Value = 5341.5
$ws.Range("M136").Value = -2791.5
$ws.Range("H137").Value = 49996.875
$ws.Range("J137").Value = 49996.875
$ws.Range("L137").Value = 49996.875
$ws.Range("N137").Value = -60196.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6739.7856
$ws.Range("I3").Value = 5719.769
$ws.Range("K3").Value = 17159.307
$ws.Range("M3").Value = -17047.307
$ws.Range("H104").Value = 299
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 62211.156
$ws.Range("J80").Value = 4327.8887
$ws.Range("L80").Value = 4327.8887
$ws.Range("N80").Value = -6323.8887
$ws.Range("H83").Value = 62211.156
$ws.Range("J83").Value = 4327.8887
$ws.Range("L83").Value = 21639.4435
$ws.Range("N83").Value = -31623.4435
$ws.Range("H102").Value = 5139.364
$ws.Range("I102").Value = 5139.364
$ws.Range("K102").Value = 5139.364
$ws.Range("M102").Value = -3517.364
$ws.Range("H122").Value = 4649.4585
$ws.Range("I122").Value = 3899.3157
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 11697.9471
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -9247.947100000001
$ws.Range("N122").Value = -27400
$ws.Range("H132").Value = 6062.42
$ws.Range("I132").Value = 5563.2563
$ws.Range("J132").Value = 7832.1816
$ws.Range("K132").Value = 16689.7689
$ws.Range("L132").Value = 23496.5448
$ws.Range("M132").Value = -14159.7689
$ws.Range("N132").Value = -28556.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 37038776
$ws.Range("I82").Value = 55557372
$ws.Range("K82").Value = 55557372
$ws.Range("M82").Value = -55557011
$ws.Range("H85").Value = 37038776
$ws.Range("I85").Value = 55557372
$ws.Range("K85").Value = 55557372
$ws.Range("M85").Value = -55556124
$ws.Range("H127").Value = 54237.617
$ws.Range("J127").Value = 54237.617
$ws.Range("L127").Value = 54237.617
$ws.Range("N127").Value = -64157.617
$ws.Range("H132").Value = 2100.7856
$ws.Range("I132").Value = 2100.7856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6302.3568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3772.3568
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2564.3333
$ws.Range("I136").Value = 2730.8462
$ws.Range("J136").Value = 2131.4
$ws.Range("K136").Value = 8192.5386
$ws.Range("L136").Value = 6394.200000000001
$ws.Range("M136").Value = -5642.5386
$ws.Range("N136").Value = -11494.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 2250
$ws.Range("I11").Value = 2000
$ws.Range("J11").Value = 2500
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = -1858
$ws.Range("N11").Value = -2784
$ws.Range("H132").Value = 4985.1304
$ws.Range("I132").Value = 3859.5454
$ws.Range("J132").Value = 7842.385
$ws.Range("K132").Value = 11578.6362
$ws.Range("L132").Value = 23527.155
$ws.Range("M132").Value = -9048.636200000001
$ws.Range("N132").Value = -28587.155
